$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 124
$wsExhibit.Range("F4").Value = 2018
$wsExhibit.Range("F5").Value = 323
$wsExhibit.Range("F6").Value = 570
$wsExhibit.Range("F8").Value = 2053
$wsExhibit.Range("F9").Value = 10434
$wsExhibit.Range("F13").Value = 200
$wsExhibit.Range("F14").Value = 403
$wsExhibit.Range("F15").Value = 7302
$wsExhibit.Range("F16").Value = 1110
$wsExhibit.Range("F17").Value = 693
$wsExhibit.Range("F18").Value = 163
$wsExhibit.Range("F19").Value = 61
$wsExhibit.Range("F20").Value = 3275

# Sheet "全部类型" (all types) - update column F ("想去人数") values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 124
$wsAll.Range("F4").Value = 2018
$wsAll.Range("F5").Value = 323
$wsAll.Range("F6").Value = 570
$wsAll.Range("F9").Value = 2053
$wsAll.Range("F12").Value = 10434
$wsAll.Range("F16").Value = 200
$wsAll.Range("F17").Value = 403
$wsAll.Range("F18").Value = 7302
$wsAll.Range("F19").Value = 1110
$wsAll.Range("F20").Value = 693
$wsAll.Range("F21").Value = 163
$wsAll.Range("F22").Value = 61
$wsAll.Range("F23").Value = 3275
